$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Semestre ideal: EQD-8,EQN-9 -> EQD-9,EQN-10
$ws.Range("B9").Value = "EQD-9,EQN-10"
$ws.Range("C9").Value = "EQD-9,EQN-10"

# Docentes responsáveis: 5840855 - Heizir Ferreira de Castro -> 1285870 - Marcos Villela Barcza
$ws.Range("B13").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C13").Value = "1285870 - Marcos Villela Barcza"

# Requisitos row 1: LOQ4002 -  Reatores Quimicos  (Requisito fraco) -> LOQ4038 -  Química Orgânica II  (Requisito fraco)
$ws.Range("B24").Value = "LOQ4038 -  Química Orgânica II  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOQ4038 -  Química Orgânica II  (Requisito fraco)`n"

# Requisitos row 2: LOT2004 -  Bioquímica  (Requisito fraco) -> LOQ4057 -  Operações Unitárias III  (Requisito fraco)
$ws.Range("B25").Value = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n"
